$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently has a "Misc Data" section (rows 15-17) containing
# "Misc Data Set#1/#2/#3". We need to insert a new "Power Factor Data"
# section (3 rows: Power Factor, Import/Export, Lagging/Leading) right
# before it, turning the old "Misc Data Set#1" row into "Lagging/Leading"
# (last row of the Power Factor section), and renumbering the remaining
# Misc Data sets down from #2/#3 to #1/#2.
# ---------------------------------------------------------------------------

# Insert two new blank rows at 15-16; everything from the old row 15 onward
# (including formulas/SUM ranges) shifts down by two rows automatically.
$ws.Rows("15:16").Insert()

# --- New row 15: Power Factor Data / Power Factor -------------------------
$ws.Range("A15").Value2 = "Power Factor Data"
$ws.Range("B15").Value2 = "Power Factor"
$ws.Range("C15").Value2 = "int16"
$ws.Range("D15").Value2 = 1
$ws.Range("E15").Value2 = 1
$ws.Range("F15").Formula = "=RIGHT(C15,2)/8*D15*E15*`$C`$1"
$ws.Range("I15").Value2 = 1000

# --- New row 16: Power Factor Data / Import/Export -------------------------
$ws.Range("A16").Value2 = "Power Factor Data"
$ws.Range("B16").Value2 = "Import/Export"
$ws.Range("C16").Value2 = "int16"
$ws.Range("D16").Value2 = 1
$ws.Range("E16").Value2 = 1
$ws.Range("F16").Formula = "=RIGHT(C16,2)/8*D16*E16*`$C`$1"
$ws.Range("I16").Value2 = 1

# --- Row 17 (was old row 15 "Misc Data Set#1"), now becomes the third
#     Power Factor Data row: Power Factor Data / Lagging/Leading ----------
$ws.Range("A17").Value2 = "Power Factor Data"
$ws.Range("B17").Value2 = "Lagging/Leading"
$ws.Range("C17").Value2 = "int16"
$ws.Range("D17").Value2 = 1
$ws.Range("E17").Value2 = 1
$ws.Range("I17").Value2 = 1
$ws.Range("H17").Clear()
$ws.Range("J17").Clear()
$ws.Range("K17").Clear()

# --- Row 18 (was old row 16 "Misc Data Set#2" with bogus DataType
#     "int17"), now becomes "Misc Data" / "Misc Data Set#1" with the
#     correct DataType -----------------------------------------------------
$ws.Range("A18").Value2 = "Misc Data"
$ws.Range("B18").Value2 = "Misc Data Set#1"
$ws.Range("C18").Value2 = "int16"

# --- Row 19 (was old row 17 "Misc Data Set#3"), renumbered to
#     "Misc Data Set#2" ------------------------------------------------------
$ws.Range("B19").Value2 = "Misc Data Set#2"

# Re-point the view selection to match the edited workbook.
$ws.Range("K28").Select()
